# Update Bordeaux neighborhood database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo and swap Latitude/Longitude column headers (C1 <-> D1)
$ws.Range("C1").Value = "Latitude"
$ws.Range("D1").Value = "Longitude"

# Update the active selection to D2
$ws.Range("D2").Select()
